$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4115
$ws.Range("I70").Value = 2985
$ws.Range("K70").Value = 8955
$ws.Range("M70").Value = -8685
$ws.Range("H73").Value = 4115
$ws.Range("I73").Value = 2985
$ws.Range("K73").Value = 8955
$ws.Range("M73").Value = -8019
$ws.Range("H103").Value = 930.5
$ws.Range("J103").Value = 1111
$ws.Range("L103").Value = 3333
$ws.Range("N103").Value = -4505
$ws.Range("H107").Value = 281.66666
$ws.Range("I107").Value = 281.66666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 281.66666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1638.33334
$ws.Range("N107").Value = $null
$ws.Range("H116").Value = 4017
$ws.Range("J116").Value = 4072.5
$ws.Range("L116").Value = 4072.5
$ws.Range("N116").Value = -10956.5
$ws.Range("H128").Value = 59750
$ws.Range("J128").Value = 59750
$ws.Range("L128").Value = 59750
$ws.Range("N128").Value = -69710
$ws.Range("H137").Value = 2519.5
$ws.Range("I137").Value = 1100
$ws.Range("J137").Value = 2992.6667
$ws.Range("K137").Value = 3300
$ws.Range("L137").Value = 8978.000100000001
$ws.Range("M137").Value = -750
$ws.Range("N137").Value = -14078.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 598.3333
$ws.Range("I2").Value = 598.3333
$ws.Range("K2").Value = 598.3333
$ws.Range("M2").Value = -485.3333
$ws.Range("H32").Value = 2145.1562
$ws.Range("I32").Value = 1891.7742
$ws.Range("K32").Value = 1891.7742
$ws.Range("M32").Value = -1604.7742
$ws.Range("H110").Value = 4499.6
$ws.Range("I110").Value = 3749.5
$ws.Range("K110").Value = 3749.5
$ws.Range("M110").Value = -1704.5
$ws.Range("H116").Value = 598.3333
$ws.Range("I116").Value = 598.3333
$ws.Range("K116").Value = 598.3333
$ws.Range("M116").Value = 1695.6667
$ws.Range("H122").Value = 2994.3333
$ws.Range("I122").Value = 2994.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8982.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6532.999899999999
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 598.3333
$ws.Range("I3").Value = 598.3333
$ws.Range("K3").Value = 598.3333
$ws.Range("M3").Value = -484.3333
$ws.Range("H20").Value = 748.73334
$ws.Range("I20").Value = 575.25
$ws.Range("K20").Value = 575.25
$ws.Range("M20").Value = -328.25
$ws.Range("H64").Value = 359
$ws.Range("I64").Value = 149
$ws.Range("J64").Value = 674
$ws.Range("K64").Value = 149
$ws.Range("L64").Value = 674
$ws.Range("M64").Value = 76
$ws.Range("N64").Value = -1124
$ws.Range("H67").Value = 359
$ws.Range("I67").Value = 149
$ws.Range("J67").Value = 674
$ws.Range("K67").Value = 149
$ws.Range("L67").Value = 674
$ws.Range("M67").Value = 631
$ws.Range("N67").Value = -2234
$ws.Range("H82").Value = 22712.076
$ws.Range("H85").Value = 22712.076
$ws.Range("H94").Value = 742.7143
$ws.Range("I94").Value = 742.7143
$ws.Range("K94").Value = 742.7143
$ws.Range("M94").Value = -291.7143
$ws.Range("H105").Value = 2390.6875
$ws.Range("I105").Value = 2417.9285
$ws.Range("K105").Value = 2417.9285
$ws.Range("M105").Value = -670.9285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
$ws.Range("H99").Value = 4390.909
$ws.Range("I99").Value = 3882.6667
$ws.Range("K99").Value = 3882.6667
$ws.Range("M99").Value = -2384.6667
$ws.Range("H126").Value = 4390.909
$ws.Range("I126").Value = 3882.6667
$ws.Range("K126").Value = 11648.0001
$ws.Range("M126").Value = -9178.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 715
$ws.Range("J5").Value = 535
$ws.Range("L5").Value = 1605
$ws.Range("N5").Value = -1829
$ws.Range("H75").Value = 13
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null
$ws.Range("H78").Value = 13
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null
$ws.Range("H92").Value = 663.8570999999999
$ws.Range("I92").Value = 629.4
$ws.Range("J92").Value = 750
$ws.Range("K92").Value = 1888.2
$ws.Range("L92").Value = 2250
$ws.Range("M92").Value = -640.1999999999998
$ws.Range("N92").Value = -4746
$ws.Range("H109").Value = 904.5
$ws.Range("I109").Value = 685.4
$ws.Range("K109").Value = 2056.2
$ws.Range("M109").Value = -1016.2
$ws.Range("H128").Value = 323940.28
$ws.Range("I128").Value = 323940.28
$ws.Range("K128").Value = 971820.8400000001
$ws.Range("M128").Value = -966840.8400000001
$ws.Range("H132").Value = 4433.3335
$ws.Range("I132").Value = 3900
$ws.Range("J132").Value = 4540
$ws.Range("K132").Value = 35100
$ws.Range("L132").Value = 40860
$ws.Range("M132").Value = -32570
$ws.Range("N132").Value = -45920
$ws.Range("H133").Value = 2000
$ws.Range("I133").Value = 2000
$ws.Range("K133").Value = 6000
$ws.Range("M133").Value = -940
$ws.Range("H135").Value = 715
$ws.Range("J135").Value = 535
$ws.Range("L135").Value = 4815
$ws.Range("N135").Value = -9885
$ws.Range("H137").Value = 2516.5
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 3033
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 9099
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -19299

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15499.7
$ws.Range("J46").Value = 21666.666
$ws.Range("L46").Value = 21666.666
$ws.Range("N46").Value = -21978.666
$ws.Range("H57").Value = 19998.334
$ws.Range("J57").Value = 19998
$ws.Range("L57").Value = 19998
$ws.Range("N57").Value = -21638
$ws.Range("H80").Value = 2394.6667
$ws.Range("I80").Value = 2308.6
$ws.Range("K80").Value = 2308.6
$ws.Range("M80").Value = -1310.6
$ws.Range("H83").Value = 2394.6667
$ws.Range("I83").Value = 2308.6
$ws.Range("K83").Value = 11543
$ws.Range("M83").Value = -6551
$ws.Range("H113").Value = 982.6667
$ws.Range("I113").Value = 982.6667
$ws.Range("K113").Value = 982.6667
$ws.Range("M113").Value = 1187.3333
$ws.Range("H126").Value = 15582
$ws.Range("J126").Value = 19999.5
$ws.Range("L126").Value = 59998.5
$ws.Range("N126").Value = -64938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1989.2858
$ws.Range("I61").Value = 1989.2858
$ws.Range("K61").Value = 1989.2858
$ws.Range("M61").Value = -1787.2858
$ws.Range("H68").Value = 2450
$ws.Range("I68").Value = 2450
$ws.Range("K68").Value = 2450
$ws.Range("M68").Value = -1701
$ws.Range("H71").Value = 2450
$ws.Range("I71").Value = 2450
$ws.Range("K71").Value = 12250
$ws.Range("M71").Value = -8506
$ws.Range("H82").Value = 1113.3334
$ws.Range("I82").Value = 1071.5
$ws.Range("J82").Value = 1197
$ws.Range("K82").Value = 1071.5
$ws.Range("L82").Value = 1197
$ws.Range("M82").Value = -710.5
$ws.Range("N82").Value = -1919
$ws.Range("H85").Value = 1113.3334
$ws.Range("I85").Value = 1071.5
$ws.Range("J85").Value = 1197
$ws.Range("K85").Value = 1071.5
$ws.Range("L85").Value = 1197
$ws.Range("M85").Value = 176.5
$ws.Range("N85").Value = -3693
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = $null
$ws.Range("H113").Value = 1989.2858
$ws.Range("I113").Value = 1989.2858
$ws.Range("K113").Value = 1989.2858
$ws.Range("M113").Value = 180.7141999999999
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = $null
